$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "This Apprach use the usual binary seach. `nFirst check if the target < mat[0][0] or target > mat[end][end] if true return false`nFirst you search the rows to see in which row the target lies, do this using binary search! `nThen after finding the row use normal binary search on the potential row in which the target lies, if found return true`nafter all this if the func is still going that means target does not exist in the mat, return false"

# Copy formatting from an existing similarly-styled row (row 3: A=plain, B=yellow fill, C=plain)
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in values in the same order the shared strings were originally appended
$ws.Range("C23").Value = $desc
$ws.Range("B23").Value = "74. Search a 2D Matrix"
$ws.Range("A23").Value = "Binary Search"

$ws.Rows("23:23").RowHeight = 86.4

[void]$ws.Range("D23").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
